$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Electric field gradient value: 4 -> 20
$ws.Range("E3").Value = 20

# Unit label changes: kV/m -> MV/m, kV -> MV
$ws.Range("F3").Value = "MV/m"
$ws.Range("F13").Value = "MV"

# Fix transit time factor formula (E12)
$ws.Range("E12").Formula = "=(2*H7)/(E9*E10)*SIN(E9*E10/2/H7)"

# E37:E41 hold literal (non-formula) comparison values that mirror H37:H41;
# update them to the recalculated results for the corrected transit-time factor.
$ws.Range("E37").Value = 0.09999785755774843
$ws.Range("E38").Value = -0.33046470148141394
$ws.Range("E39").Value = -0.1999986695317941
$ws.Range("E40").Value = 0.8070662800280045
$ws.Range("E41").Value = 0.2000034938176917

# Re-touch the K37:K41 "Difference" formulas so they pick up the refreshed
# E37:E41 / H37:H41 values (they stay functionally identical: =E{n}-H{n}).
$ws.Range("K37").Formula = "=E37-H37"
$ws.Range("K38").Formula = "=E38-H38"
$ws.Range("K39").Formula = "=E39-H39"
$ws.Range("K40").Formula = "=E40-H40"
$ws.Range("K41").Formula = "=E41-H41"

# Update the cell selection to match the saved state of the workbook
$ws.Range("C47").Select()
